$wb = $excel.ActiveWorkbook

# Sheet 1: "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 220
$ws1.Range("F3").Value = 1047
$ws1.Range("F4").Value = 521
$ws1.Range("F5").Value = 13768
$ws1.Range("F6").Value = 43
$ws1.Range("F9").Value = 162
$ws1.Range("F13").Value = 519
$ws1.Range("F15").Value = 13793
$ws1.Range("F16").Value = 357
$ws1.Range("F17").Value = 616
$ws1.Range("F18").Value = 9081
$ws1.Range("F20").Value = 8196
$ws1.Range("F21").Value = 267
$ws1.Range("F35").Value = 207
$ws1.Range("F38").Value = 5055

# Sheet 2: "演出" (performances) - update column F
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 47

# Sheet 4: "全部类型" (all types) - update column F
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 220
$ws4.Range("F3").Value = 1047
$ws4.Range("F4").Value = 521
$ws4.Range("F5").Value = 13768
$ws4.Range("F6").Value = 43
$ws4.Range("F9").Value = 162
$ws4.Range("F13").Value = 519
$ws4.Range("F15").Value = 13793
$ws4.Range("F16").Value = 357
$ws4.Range("F17").Value = 616
$ws4.Range("F18").Value = 9081
$ws4.Range("F20").Value = 8196
$ws4.Range("F21").Value = 267
$ws4.Range("F32").Value = 47
$ws4.Range("F37").Value = 207
$ws4.Range("F40").Value = 5055
